# Applies:
#  1) Date placeholder bump 29/07/2025 -> 31/07/2025 on the slide master
#     and every slide layout (datetimeFigureOut field caches the date as
#     literal text; PowerPoint recalculates/re-caches it whenever the
#     placeholder's text is touched).
#  2) "Endosymbiote" -> "Endosymbiont" spelling fix on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "29/07/2025"
$newDate = "31/07/2025"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        if ($shp.Name -notlike "Date Placeholder*") { continue }

        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $chars = $tr.Characters(1, $tr.Length)
            $chars.Text = $newDate
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom (slide) layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Slide 1: fix "Endosymbiote" -> "Endosymbiont" inside the grouped figure text.
$slide = $p.Slides.Item(1)

function Fix-Endosymbiote {
    param($shape)

    if ($shape.Type -eq 6) {
        # msoGroup - recurse into its members.
        $items = $shape.GroupItems
        for ($j = 1; $j -le $items.Count; $j++) {
            Fix-Endosymbiote $items.Item($j)
        }
        return
    }

    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -match "Endosymbiote") {
        $tr.Replace("Endosymbiote", "Endosymbiont", 0, 0, 0)
    }
}

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    Fix-Endosymbiote $slide.Shapes.Item($i)
}
